# Phase0/CRC Model.pptx - slide 6 ("Check" CRC card group):
#   - "Use the Items in Player's inventory" bullet -> "Check against Player attributes"
#   - "Add or discard any Item" bullet and the trailing empty bullet paragraph removed
#   - "Create new Check" bullet text is unchanged
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)

# The three bullet lines live in shape "Google Shape;59;p13", inside the
# "Group 8" group shape that sits on this slide.
$group = $s.Shapes.Item("Group 8")
$card = $group.GroupItems.Item("Google Shape;59;p13")

# Replace the 4-paragraph bullet list (incl. the blank trailing paragraph)
# with the new 2-paragraph list in one shot.
$card.TextFrame.TextRange.Text = "Create new Check`rCheck against Player attributes"
